$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.167.85'
$ws.Range("E2").Value = '  +0.34%  '

# Row 3
$ws.Range("D3").Value = '1.909.31'
$ws.Range("E3").Value = '  -0.12%  '

# Row 4
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.8208'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +4.47%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '243.75'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.33%  '

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3258'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +2.95%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '26.84'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +2.29%  '

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.07048'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +1.75%  '

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.08100'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +1.48%  '

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.7709'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +3.37%  '

# Row 13
$ws.Range("D13").Value = '1.915.02'
$ws.Range("E13").Value = '  +0.15%  '

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '5.291'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +1.28%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '92.97'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -0.39%  '

# Row 16
$ws.Range("D16").Value = '30.174.78'
$ws.Range("E16").Value = '  +0.29%  '

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '14.22'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +1.49%  '

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '5.899'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -0.33%  '

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '245.68'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -0.58%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.000007791'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.16%  '

# Row 21
$ws.Range("D21").Value = '2.164.06'
$ws.Range("E21").Value = '  -0.01%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.17%  '

# Row 23
$ws.Range("E23").Value = '  +0.08%  '

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '7.046'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +1.90%  '

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.1662'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +20.67%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '9.323'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.18%  '

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '167.38'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -1.21%  '

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '18.98'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +0.35%  '

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '2.109'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +3.16%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.371'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.31%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.527'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.26%  '

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.05913'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +2.74%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '4.305'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -0.78%  '

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '4.100'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -0.42%  '

# Row 35
$ws.Range("E35").Value = '  +0.93%  '

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.7345'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.31%  '

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.715'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.28%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.01928'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.28%  '

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.798'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +0.15%  '

# Row 40
$ws.Range("E40").Value = '  +0.58%  '

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '73.36'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +1.06%  '

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '5.963'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -3.30%  '

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.8530'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +2.15%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '1.909'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.62%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.07%  '

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '102.73'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +2.14%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '9.877'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.60%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '7.594'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.38%  '

# Row 49
$ws.Range("D49").Value = '1.007.62'
$ws.Range("E49").Value = '  +1.76%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '1.573'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +5.21%  '

# Row 51
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.061.23'
$ws.Range("E51").Value = '  +0.23%  '
